$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" between "总计" and "2021-Q4".
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q3"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Re-fetch sheet handles by name now that the sheet collection has shifted.
$totalSheet = $wb.Worksheets.Item("总计")
$oldSheet = $wb.Worksheets.Item("2021-Q4")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# 2. "总计" sheet: push the existing 2021-Q4 summary row down to row 3 and
#    write the new 2022-Q3 summary row into row 2.
# ---------------------------------------------------------------------------

# Clone the formatting of A2 (the styled index cell) down into A3 first, so
# the moved row keeps the same look the original row had.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 9
$totalSheet.Range("D3").Value = 1.62

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.19

# ---------------------------------------------------------------------------
# 3. Build the new "2022-Q3" sheet content, reusing the existing "2021-Q4"
#    sheet's header/index-column formatting.
# ---------------------------------------------------------------------------

$oldSheet.Range("B1:H1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$oldSheet.Range("A2:A6").Copy()
$q3Sheet.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
$q3Sheet.Range("A2").Value = 0
Set-TextValue $q3Sheet.Range("B2") "009225"
$q3Sheet.Range("C2").Value = "天弘中证中美互联网指数（QDII）A"
Set-TextValue $q3Sheet.Range("D2") "1.20"
Set-TextValue $q3Sheet.Range("E2") "94.98"
Set-TextValue $q3Sheet.Range("F2") "7.98"
Set-TextValue $q3Sheet.Range("G2") "0.0958"
$q3Sheet.Range("H2").Value = 5

# Row 3
$q3Sheet.Range("A3").Value = 1
Set-TextValue $q3Sheet.Range("B3") "009226"
$q3Sheet.Range("C3").Value = "天弘中证中美互联网指数（QDII）C"
Set-TextValue $q3Sheet.Range("D3") "0.60"
Set-TextValue $q3Sheet.Range("E3") "94.98"
Set-TextValue $q3Sheet.Range("F3") "7.98"
Set-TextValue $q3Sheet.Range("G3") "0.0479"
$q3Sheet.Range("H3").Value = 5

# Row 4
$q3Sheet.Range("A4").Value = 2
Set-TextValue $q3Sheet.Range("B4") "012751"
$q3Sheet.Range("C4").Value = "建信纳斯达克100指数（QDII）A 美元现汇"
Set-TextValue $q3Sheet.Range("D4") "0.64"
Set-TextValue $q3Sheet.Range("E4") "80.13"
Set-TextValue $q3Sheet.Range("F4") "2.37"
Set-TextValue $q3Sheet.Range("G4") "0.0152"
$q3Sheet.Range("H4").Value = 8

# Row 5
$q3Sheet.Range("A5").Value = 3
Set-TextValue $q3Sheet.Range("B5") "012752"
$q3Sheet.Range("C5").Value = "建信纳斯达克100指数（QDII）C 人民币"
Set-TextValue $q3Sheet.Range("D5") "0.64"
Set-TextValue $q3Sheet.Range("E5") "80.13"
Set-TextValue $q3Sheet.Range("F5") "2.37"
Set-TextValue $q3Sheet.Range("G5") "0.0152"
$q3Sheet.Range("H5").Value = 8

# Row 6
$q3Sheet.Range("A6").Value = 4
Set-TextValue $q3Sheet.Range("B6") "012753"
$q3Sheet.Range("C6").Value = "建信纳斯达克100指数（QDII）C 美元现汇"
Set-TextValue $q3Sheet.Range("D6") "0.64"
Set-TextValue $q3Sheet.Range("E6") "80.13"
Set-TextValue $q3Sheet.Range("F6") "2.37"
Set-TextValue $q3Sheet.Range("G6") "0.0152"
$q3Sheet.Range("H6").Value = 8
